$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 122, shifting existing rows 122+ down by one.
$ws.Rows.Item(122).Insert()

# Populate the new row 122 with its values (same shape as neighboring rows).
$ws.Cells.Item(122, 1).Value = 10
$ws.Cells.Item(122, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(122, 3).Value = "La Araucanía"
$ws.Cells.Item(122, 4).Value = 44518
$ws.Cells.Item(122, 4).NumberFormat = $ws.Cells.Item(121, 4).NumberFormat
$ws.Cells.Item(122, 5).Value = 9
$ws.Cells.Item(122, 6).Value = 100112039
$ws.Cells.Item(122, 7).Value = "Ciboulette"
$ws.Cells.Item(122, 8).Value = "Sin especificar"
$ws.Cells.Item(122, 9).Value = "Primera"
$ws.Cells.Item(122, 10).Value = 40
$ws.Cells.Item(122, 11).Value = 5000
$ws.Cells.Item(122, 12).Value = 5000
$ws.Cells.Item(122, 13).Value = 5000
$ws.Cells.Item(122, 14).Value = "`$/docena de atados"
$ws.Cells.Item(122, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(122, 16).Value = 1667
$ws.Cells.Item(122, 17).Value = 3
$ws.Cells.Item(122, 18).Value = "Hortaliza"
